$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "B"=18.73820516028269; "C"=5.745154400776977; "D"=8.209324402614632; "F"=41.05517398452946; "G"=3.726632471681413; "K"=14.68741086150565; "L"=10.91131242402436; "N"=23.44912315802864 }
  3 = @{ "B"=18.53954331618607; "C"=5.641558213594811; "D"=8.220266429876469; "F"=40.88741461769141; "G"=3.729835017758365; "K"=14.55737814277456; "L"=10.89162186166533; "N"=23.47735831792736 }
  4 = @{ "B"=18.42188522008228; "C"=5.57589823737071; "D"=8.227174631837373; "F"=40.79266645356797; "G"=3.731903849753103; "K"=14.48106123731799; "L"=10.88181371299923; "N"=23.49651954655245 }
  5 = @{ "B"=18.37507896468263; "C"=5.548639000161221; "D"=8.230037737192097; "F"=40.75615402430899; "G"=3.732772771793859; "K"=14.45087949460958; "L"=10.87839300601252; "N"=23.50478644137936 }
  6 = @{ "B"=18.3673772207608; "C"=5.544082699659052; "D"=8.230516057210261; "F"=40.75021848029454; "G"=3.732918620016538; "K"=14.44592418170364; "L"=10.87785986121472; "N"=23.50618683856235 }
  7 = @{ "B"=18.42124928801841; "C"=5.575532624927595; "D"=8.227213050202369; "F"=40.79216551167266; "G"=3.73191546352658; "K"=14.48065043798824; "L"=10.88176524436899; "N"=23.49662918065343 }
  8 = @{ "B"=18.66884272967786; "C"=5.709868168578104; "D"=8.213057977833751; "F"=40.99563084134433; "G"=3.72771550088629; "K"=14.64186396497997; "L"=10.90405106892609; "N"=23.45847973742581 }
  9 = @{ "B"=19.18618435575958; "C"=5.956462061966145; "D"=8.186793918473374; "F"=41.45905761860767; "G"=3.720288090491606; "K"=14.98447282718697; "L"=10.96573544831491; "N"=23.39816136482858 }
  10 = @{ "B"=19.58223273768274; "C"=6.126674010610166; "D"=8.168391431333873; "F"=41.83720673115997; "G"=3.715318249452951; "K"=15.25022580533348; "L"=11.02182265185929; "N"=23.36270083490989 }
  11 = @{ "B"=19.76511756436103; "C"=6.201595103654616; "D"=8.160210015799523; "F"=42.01702374312722; "G"=3.713161839362581; "K"=15.37370252336777; "L"=11.04962519895157; "N"=23.34849615447066 }
  12 = @{ "B"=19.83469888434831; "C"=6.229594704854659; "D"=8.157138972526127; "F"=42.08620005179786; "G"=3.712360178159141; "K"=15.42079117964347; "L"=11.06047725798473; "N"=23.34339458904645 }
  13 = @{ "B"=19.81969978964265; "C"=6.223581168102838; "D"=8.157799175947117; "F"=42.07125412321601; "G"=3.71253216786001; "K"=15.41063574167523; "L"=11.05812575762011; "N"=23.344480956782 }
  14 = @{ "B"=19.7708358796479; "C"=6.203906148900031; "D"=8.159956817825147; "F"=42.0226934258495; "G"=3.713095587603292; "K"=15.37757014318436; "L"=11.05051154554183; "N"=23.3480708829988 }
  15 = @{ "B"=19.74094597674416; "C"=6.191805964098209; "D"=8.161281955549349; "F"=41.99308855723273; "G"=3.713442639487114; "K"=15.35735835661123; "L"=11.04588963214596; "N"=23.35030595926289 }
  16 = @{ "B"=19.57033001776807; "C"=6.121726329386091; "D"=8.168929908547641; "F"=41.82560897094283; "G"=3.71546126914228; "K"=15.24220482511934; "L"=11.02005129980108; "N"=23.36366794547135 }
  17 = @{ "B"=19.46631356496338; "C"=6.078084328037646; "D"=8.173670167415599; "F"=41.72483619916292; "G"=3.716726308873498; "K"=15.17219478499251; "L"=11.00478282167815; "N"=23.37235880711099 }
  18 = @{ "B"=19.4067459687599; "C"=6.05274701328506; "D"=8.176414533745293; "F"=41.66761044893127; "G"=3.7174637576693; "K"=15.13217259463295; "L"=10.99621651283099; "N"=23.37753886735289 }
  19 = @{ "B"=19.38662400956401; "C"=6.044128105247174; "D"=8.177346809663243; "F"=41.64836232301263; "G"=3.717715136384519; "K"=15.11866521430566; "L"=10.99335330849103; "N"=23.37932387361545 }
  20 = @{ "B"=19.47735986073889; "C"=6.082754555084134; "D"=8.1731637084672; "F"=41.73548771595047; "G"=3.716590626359448; "K"=15.17962233626046; "L"=11.00638588775874; "N"=23.37141488263744 }
  21 = @{ "B"=19.78518003693358; "C"=6.209695339072065; "D"=8.159322333004903; "F"=42.03692778159614; "G"=3.712929693215165; "K"=15.38727366347059; "L"=11.05273927896665; "N"=23.34700890179149 }
  22 = @{ "B"=19.98823359444532; "C"=6.290489336304926; "D"=8.150433934779103; "F"=42.24023310528696; "G"=3.710624012917434; "K"=15.52489425583058; "L"=11.08491875594793; "N"=23.33267549378715 }
  23 = @{ "B"=19.87970968573459; "C"=6.247569842012543; "D"=8.155163484843559; "F"=42.13116192869425; "G"=3.71184667062434; "K"=15.45128251971737; "L"=11.06757332996164; "N"=23.34017738530985 }
  24 = @{ "B"=19.47236510022248; "C"=6.080643913843383; "D"=8.173392618941339; "F"=41.73066995120331; "G"=3.716651936768087; "K"=15.17626362854267; "L"=11.00566048165589; "N"=23.37184105908701 }
  25 = @{ "B"=19.04317373060362; "C"=5.891629103947133; "D"=8.193740864555108; "F"=41.32694696616655; "G"=3.722211439296397; "K"=14.88916606399816; "L"=10.94714130711607; "N"=23.41292543233441 }
}

foreach ($rowKey in $data.Keys) {
  $rowData = $data[$rowKey]
  foreach ($col in $rowData.Keys) {
    $ws.Range("$col$rowKey").Value = $rowData[$col]
  }
}
